$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5237.875
$ws.Range("I19").Value = 5475
$ws.Range("K19").Value = 5475
$ws.Range("M19").Value = -5300
$ws.Range("H52").Value = 4009
$ws.Range("I52").Value = 4009
$ws.Range("K52").Value = 12027
$ws.Range("M52").Value = -11867
$ws.Range("H62").Value = 4198.231
$ws.Range("I62").Value = 3035.4
$ws.Range("K62").Value = 3035.4
$ws.Range("M62").Value = -2411.4
$ws.Range("H65").Value = 4198.231
$ws.Range("I65").Value = 3035.4
$ws.Range("K65").Value = 15177
$ws.Range("M65").Value = -12057
$ws.Range("H74").Value = 4235.8667
$ws.Range("I74").Value = 3609.375
$ws.Range("K74").Value = 3609.375
$ws.Range("M74").Value = -2673.375
$ws.Range("H77").Value = 4235.8667
$ws.Range("I77").Value = 3609.375
$ws.Range("K77").Value = 18046.875
$ws.Range("M77").Value = -13366.875
$ws.Range("H100").Value = 2314.0527
$ws.Range("I100").Value = 1281.5385
$ws.Range("J100").Value = 4551.1665
$ws.Range("K100").Value = 1281.5385
$ws.Range("L100").Value = 4551.1665
$ws.Range("M100").Value = -740.5385000000001
$ws.Range("N100").Value = -5633.1665
$ws.Range("H116").Value = 117959.055
$ws.Range("I116").Value = 147609.58
$ws.Range("J116").Value = 99090.55
$ws.Range("K116").Value = 147609.58
$ws.Range("L116").Value = 99090.55
$ws.Range("M116").Value = -144167.58
$ws.Range("N116").Value = -105974.55
$ws.Range("H138").Value = 3640.4324
$ws.Range("I138").Value = 1434.5454
$ws.Range("K138").Value = 4303.6362
$ws.Range("M138").Value = 836.3638000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1287.7632
$ws.Range("I32").Value = 1342.4166
$ws.Range("K32").Value = 1342.4166
$ws.Range("M32").Value = -1055.4166
$ws.Range("H45").Value = 1239.1666
$ws.Range("I45").Value = 1239.1666
$ws.Range("K45").Value = 1239.1666
$ws.Range("M45").Value = -862.1666
$ws.Range("H61").Value = 2867.3
$ws.Range("I61").Value = 2517.6667
$ws.Range("K61").Value = 2517.6667
$ws.Range("M61").Value = -2305.6667
$ws.Range("H74").Value = 4210234.5
$ws.Range("I74").Value = 3704661.2
$ws.Range("J74").Value = 4631546
$ws.Range("K74").Value = 3704661.2
$ws.Range("L74").Value = 4631546
$ws.Range("M74").Value = -3703787.2
$ws.Range("N74").Value = -4633294
$ws.Range("H77").Value = 4210234.5
$ws.Range("I77").Value = 3704661.2
$ws.Range("J77").Value = 4631546
$ws.Range("K77").Value = 18523306
$ws.Range("L77").Value = 23157730
$ws.Range("M77").Value = -18518938
$ws.Range("N77").Value = -23166466
$ws.Range("H102").Value = 1457.5625
$ws.Range("I102").Value = 1457.5625
$ws.Range("K102").Value = 1457.5625
$ws.Range("M102").Value = 164.4375
$ws.Range("H122").Value = 3851.889
$ws.Range("I122").Value = 2111.3333
$ws.Range("J122").Value = 7333
$ws.Range("K122").Value = 6333.999899999999
$ws.Range("L122").Value = 21999
$ws.Range("M122").Value = -3883.999899999999
$ws.Range("N122").Value = -26899
$ws.Range("H124").Value = 43450
$ws.Range("J124").Value = 43450
$ws.Range("L124").Value = 43450
$ws.Range("N124").Value = -53270
$ws.Range("H136").Value = 2867.3
$ws.Range("I136").Value = 2517.6667
$ws.Range("K136").Value = 7553.000100000001
$ws.Range("M136").Value = -5003.000100000001
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").Value = ""
$ws.Range("H36").Value = 7406.4287
$ws.Range("I36").Value = 2352.6
$ws.Range("K36").Value = 2352.6
$ws.Range("M36").Value = -1818.6
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = ""
$ws.Range("H104").Value = 45000
$ws.Range("J104").Value = 45000
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -51988
$ws.Range("H105").Value = 3269.077
$ws.Range("I105").Value = 3959.4
$ws.Range("J105").Value = 2837.625
$ws.Range("K105").Value = 3959.4
$ws.Range("L105").Value = 2837.625
$ws.Range("M105").Value = -2212.4
$ws.Range("N105").Value = -6331.625
$ws.Range("H124").Value = 52000
$ws.Range("J124").Value = 52000
$ws.Range("L124").Value = 52000
$ws.Range("N124").Value = -61820
$ws.Range("H134").Value = 15354227
$ws.Range("I134").Value = 8067850.5
$ws.Range("J134").Value = 47622464
$ws.Range("K134").Value = 24203551.5
$ws.Range("L134").Value = 142867392
$ws.Range("M134").Value = -24201016.5
$ws.Range("N134").Value = -142872462

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 38665
$ws.Range("J28").Value = 38665
$ws.Range("L28").Value = 38665
$ws.Range("N28").Value = -39155
$ws.Range("H31").Value = 13675.044
$ws.Range("I31").Value = 22128.166
$ws.Range("J31").Value = 4453.4546
$ws.Range("K31").Value = 22128.166
$ws.Range("L31").Value = 4453.4546
$ws.Range("M31").Value = -21833.166
$ws.Range("N31").Value = -5043.4546
$ws.Range("H34").Value = 13675.044
$ws.Range("I34").Value = 22128.166
$ws.Range("J34").Value = 4453.4546
$ws.Range("K34").Value = 22128.166
$ws.Range("L34").Value = 4453.4546
$ws.Range("M34").Value = -21926.166
$ws.Range("N34").Value = -4857.4546
$ws.Range("H105").Value = 4011.8462
$ws.Range("I105").Value = 4312
$ws.Range("K105").Value = 4312
$ws.Range("M105").Value = -2565
$ws.Range("H138").Value = 98689.75
$ws.Range("J138").Value = 98689.75
$ws.Range("L138").Value = 98689.75
$ws.Range("N138").Value = -108969.75

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 659.7
$ws.Range("J2").Value = 165.64285
$ws.Range("L2").Value = 993.8571000000001
$ws.Range("N2").Value = -1219.8571
$ws.Range("H54").Value = 47999.8
$ws.Range("J54").Value = 47999.8
$ws.Range("L54").Value = 143999.4
$ws.Range("N54").Value = -145117.4
$ws.Range("H86").Value = 320.14285
$ws.Range("J86").Value = 330.6
$ws.Range("L86").Value = 991.8000000000001
$ws.Range("N86").Value = -3363.8
$ws.Range("H89").Value = 320.14285
$ws.Range("J89").Value = 330.6
$ws.Range("L89").Value = 2975.4
$ws.Range("N89").Value = -14831.4
$ws.Range("H131").Value = 628251
$ws.Range("J131").Value = 910515.9
$ws.Range("L131").Value = 2731547.7
$ws.Range("N131").Value = -2741627.7
$ws.Range("H137").Value = 5554
$ws.Range("I137").Value = 4999.5
$ws.Range("J137").Value = 5738.8335
$ws.Range("K137").Value = 14998.5
$ws.Range("L137").Value = 17216.5005
$ws.Range("M137").Value = -9898.5
$ws.Range("N137").Value = -27416.5005
$ws.Range("H139").Value = 2966.2
$ws.Range("I139").Value = 2014.6666
$ws.Range("K139").Value = 6043.9998
$ws.Range("M139").Value = -903.9997999999996

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3851.8647
$ws.Range("I102").Value = 3748.818
$ws.Range("K102").Value = 3748.818
$ws.Range("M102").Value = -2126.818
$ws.Range("H107").Value = 585.5833
$ws.Range("I107").Value = 313.1
$ws.Range("K107").Value = 313.1
$ws.Range("M107").Value = 1606.9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2967.3333
$ws.Range("I122").Value = 2407.6667
$ws.Range("K122").Value = 7223.000100000001
$ws.Range("M122").Value = -4773.000100000001
$ws.Range("H136").Value = 27781634
$ws.Range("I136").Value = 3048.1853
$ws.Range("K136").Value = 9144.555899999999
$ws.Range("M136").Value = -6594.555899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 70000
$ws.Range("J116").Value = 70000
$ws.Range("L116").Value = 70000
$ws.Range("N116").Value = -79178
$ws.Range("H122").Value = 3205.9524
$ws.Range("I122").Value = 3409.2
$ws.Range("J122").Value = 2697.8333
$ws.Range("K122").Value = 10227.6
$ws.Range("L122").Value = 8093.499899999999
$ws.Range("M122").Value = -7777.599999999999
$ws.Range("N122").Value = -12993.4999
